$wb = $excel.ActiveWorkbook

# 1. Change the selection on the SolarPV sheet before touching anything else,
#    so it does not end up being the "active" sheet in the saved file.
$solarPV = $wb.Worksheets.Item("SolarPV")
$solarPV.Range("D1:D2").Select() | Out-Null

# 2. Insert a new "Battery" worksheet right before "Switch" (i.e. right after
#    "SolarPV"), mirroring how Excel numbers a newly inserted sheet (next
#    free sheetId) and renumbers the relationship id used by "Switch".
$afterSheet = $wb.Worksheets.Item("SolarPV")
$battery = $wb.Worksheets.Add($null, $afterSheet)
$battery.Name = "Battery"

# 3. Populate the header row + data row. Write the "Battery1" label before
#    the other new headers so the shared-string table is built in the same
#    order as the source edit (Battery1, E0, Emax, SOCmin, SOCmax).
$battery.Range("A1").Value = "Name"
$battery.Range("A2").Value = "Battery1"
$battery.Range("B1").Value = "E0"
$battery.Range("C1").Value = "Emax"
$battery.Range("D1").Value = "SOCmin"
$battery.Range("E1").Value = "SOCmax"
$battery.Range("F1").Value = "Pmax"
$battery.Range("G1").Value = "CONNECTION"

$battery.Range("B2").Value = 7000
$battery.Range("C2").Value = 15000
$battery.Range("D2").Value = 0.1
$battery.Range("E2").Value = 0.98
$battery.Range("F2").Value = 1000
$battery.Range("G2").Value = "P,EB1,P;"

# 4. Leave the cursor where the author left it and make "Battery" the active
#    (selected) sheet/tab, matching the new tabSelected + activeTab state.
$battery.Range("F10").Select() | Out-Null
$battery.Activate() | Out-Null
